# Apply updated cryptocurrency price/volume data to Sheet1
# (values that look numeric must be forced to Text format so they
#  are preserved exactly as strings, matching the source data feed)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '61.714.09'
$ws.Range('D3').Value = '3.402.63'
$ws.Range('E3').Value = '  -0.19%  '
$ws.Range('E4').Value = '  -0.11%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '412.26'
$ws.Range('E5').Value = '  +1.00%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '129.37'
$ws.Range('E6').Value = '  +0.89%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.619'
$ws.Range('E7').Value = '  -2.02%  '
$ws.Range('E8').Value = '  +0.03%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.721'
$ws.Range('E9').Value = '  -1.49%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.135'
$ws.Range('E10').Value = '  -4.54%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '42.51'
$ws.Range('E11').Value = '  +0.59%  '
$ws.Range('B12').Value = 'Polkadot'
$ws.Range('C12').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '9.11'
$ws.Range('E12').Value = '  +2.37%  '
$ws.Range('B13').Value = 'ShibaInu'
$ws.Range('C13').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.0000215'
$ws.Range('E13').Value = '  -1.54%  '
$ws.Range('D14').Value = '3.948.34'
$ws.Range('E14').Value = '  -0.19%  '
$ws.Range('E15').Value = '  -0.05%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '20.35'
$ws.Range('E16').Value = '  -1.62%  '
$ws.Range('D17').Value = '3.405.70'
$ws.Range('E17').Value = '  -0.02%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '12.41'
$ws.Range('E18').Value = '  +2.91%  '
$ws.Range('E19').Value = '  +0.47%  '
$ws.Range('D20').Value = '61.769.40'
$ws.Range('E20').Value = '  -0.17%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '480.76'
$ws.Range('E21').Value = '  +16.73%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '90.66'
$ws.Range('E22').Value = '  +1.84%  '
$ws.Range('E23').Value = '  +3.09%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '13.05'
$ws.Range('E24').Value = '  +0.33%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '3.31'
$ws.Range('E25').Value = '  +2.51%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '9.75'
$ws.Range('E26').Value = '  +10.39%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '32.97'
$ws.Range('E27').Value = '  -0.12%  '
$ws.Range('E28').Value = '  -0.56%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '7.71'
$ws.Range('E29').Value = '  +1.82%  '
$ws.Range('E30').Value = '  -2.70%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '11.81'
$ws.Range('E31').Value = '  -0.33%  '
$ws.Range('E32').Value = '  -2.03%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.111'
$ws.Range('E33').Value = '  -3.18%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '40.91'
$ws.Range('E34').Value = '  -4.01%  '
$ws.Range('E35').Value = '  -0.71%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '58.37'
$ws.Range('E36').Value = '  +7.98%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.0484'
$ws.Range('E37').Value = '  -2.20%  '
$ws.Range('E38').Value = '  +0.17%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '3.01'
$ws.Range('E39').Value = '  +3.89%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '148.82'
$ws.Range('E40').Value = '  +5.10%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.321'
$ws.Range('E41').Value = '  +3.44%  '
$ws.Range('E42').Value = '  +0.62%  '
$ws.Range('E43').Value = '  -0.27%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '2.05'
$ws.Range('E44').Value = '  +4.38%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '2.58'
$ws.Range('E45').Value = '  +7.07%  '
$ws.Range('E46').Value = '  +2.71%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '2.32'
$ws.Range('E47').Value = '  +18.17%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '16.38'
$ws.Range('E48').Value = '  -1.09%  '
$ws.Range('B49').Value = 'EnergySwap'
$ws.Range('C49').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '22.11'
$ws.Range('E49').Value = '  +1.78%  '
$ws.Range('B50').Value = 'PEPE'
$ws.Range('C50').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D50').Value = '0.0₃0516'
$ws.Range('E50').Value = '  +15.22%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '112.47'
$ws.Range('E51').Value = '  +14.25%  '
